# Atualizacao de bases das ligas (swap/rotate match rows 75/76, 147/148, 314/315/316)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 75
$ws.Range("B75").Value = 6803727
$arr = New-Object 'object[,]' 1,26
$arr[0,0] = 'Chrobry Glogow'
$arr[0,1] = 'GKS Tychy 71'
$arr[0,2] = 2
$arr[0,3] = 1
$arr[0,4] = 0
$arr[0,5] = 1
$arr[0,6] = 'H'
$arr[0,7] = 3.1
$arr[0,8] = 3.3
$arr[0,9] = 2.15
$arr[0,10] = 3.75
$arr[0,11] = 3.6
$arr[0,12] = 1.85
$arr[0,13] = 0.5
$arr[0,14] = 1.925
$arr[0,15] = 1.875
$arr[0,16] = 2.75
$arr[0,17] = 1.925
$arr[0,18] = 1.875
$arr[0,19] = 2.75
$arr[0,20] = -1
$arr[0,21] = -1
$arr[0,22] = 0.925
$arr[0,23] = -1
$arr[0,24] = 0.4625
$arr[0,25] = -0.5
$ws.Range("E75:AD75").Value = $arr

# Row 76
$ws.Range("B76").Value = 6805719
$arr = New-Object 'object[,]' 1,26
$arr[0,0] = 'Motor Lublin'
$arr[0,1] = 'Stal Rzeszow'
$arr[0,2] = 3
$arr[0,3] = 2
$arr[0,4] = 1
$arr[0,5] = 1
$arr[0,6] = 'H'
$arr[0,7] = 2.3
$arr[0,8] = 3.3
$arr[0,9] = 2.8
$arr[0,10] = 2.05
$arr[0,11] = 3.4
$arr[0,12] = 3.2
$arr[0,13] = -0.25
$arr[0,14] = 1.85
$arr[0,15] = 2
$arr[0,16] = 2.5
$arr[0,17] = 1.875
$arr[0,18] = 1.975
$arr[0,19] = 1.05
$arr[0,20] = -1
$arr[0,21] = -1
$arr[0,22] = 0.8500000000000001
$arr[0,23] = -1
$arr[0,24] = 0.875
$arr[0,25] = -1
$ws.Range("E76:AD76").Value = $arr

# Row 147
$ws.Range("B147").Value = 6803778
$arr = New-Object 'object[,]' 1,26
$arr[0,0] = 'Podbeskidzie Bielsko Biala'
$arr[0,1] = 'Miedz Legnica'
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 'D'
$arr[0,7] = 3.1
$arr[0,8] = 3.4
$arr[0,9] = 2.1
$arr[0,10] = 3
$arr[0,11] = 3.4
$arr[0,12] = 2.15
$arr[0,13] = 0.25
$arr[0,14] = 1.9
$arr[0,15] = 1.95
$arr[0,16] = 2.5
$arr[0,17] = 1.85
$arr[0,18] = 2
$arr[0,19] = -1
$arr[0,20] = 2.4
$arr[0,21] = -1
$arr[0,22] = 0.45
$arr[0,23] = -0.5
$arr[0,24] = -1
$arr[0,25] = 1
$ws.Range("E147:AD147").Value = $arr

# Row 148
$ws.Range("B148").Value = 6803779
$arr = New-Object 'object[,]' 1,26
$arr[0,0] = 'Zaglebie Sosnowiec'
$arr[0,1] = 'Arka Gdynia'
$arr[0,2] = 1
$arr[0,3] = 3
$arr[0,4] = 1
$arr[0,5] = 1
$arr[0,6] = 'A'
$arr[0,7] = 3.3
$arr[0,8] = 3.4
$arr[0,9] = 2
$arr[0,10] = 4.2
$arr[0,11] = 3.5
$arr[0,12] = 1.727
$arr[0,13] = 0.75
$arr[0,14] = 1.825
$arr[0,15] = 2.025
$arr[0,16] = 2.5
$arr[0,17] = 1.9
$arr[0,18] = 1.95
$arr[0,19] = -1
$arr[0,20] = -1
$arr[0,21] = 0.7270000000000001
$arr[0,22] = -1
$arr[0,23] = 1.025
$arr[0,24] = 0.8999999999999999
$arr[0,25] = -1
$ws.Range("E148:AD148").Value = $arr

# Row 314
$ws.Range("B314").Value = 7096878
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 'Miedz Legnica'
$arr[0,1] = 'Lechia Gdansk'
$arr[0,2] = 4
$arr[0,3] = 1
$ws.Range("E314:H314").Value = $arr
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = 'H'
$arr[0,1] = 2.5
$arr[0,2] = 3.2
$arr[0,3] = 2.5
$arr[0,4] = 1.909
$arr[0,5] = 3.25
$arr[0,6] = 3.6
$arr[0,7] = -0.5
$arr[0,8] = 2
$arr[0,9] = 1.85
$arr[0,10] = 2.75
$arr[0,11] = 1.9
$arr[0,12] = 1.95
$arr[0,13] = 0.909
$arr[0,14] = -1
$arr[0,15] = -1
$arr[0,16] = 1
$arr[0,17] = -1
$arr[0,18] = 0.8999999999999999
$arr[0,19] = -1
$ws.Range("K314:AD314").Value = $arr

# Row 315
$ws.Range("B315").Value = 7093053
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 'Stal Rzeszow'
$arr[0,1] = 'Polonia Warsaw'
$arr[0,2] = 1
$arr[0,3] = 2
$ws.Range("E315:H315").Value = $arr
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = 'A'
$arr[0,1] = 4
$arr[0,2] = 4.333
$arr[0,3] = 1.571
$arr[0,4] = 3.9
$arr[0,5] = 4.5
$arr[0,6] = 1.571
$arr[0,7] = 1
$arr[0,8] = 1.775
$arr[0,9] = 2.025
$arr[0,10] = 3.25
$arr[0,11] = 1.9
$arr[0,12] = 1.9
$arr[0,13] = -1
$arr[0,14] = -1
$arr[0,15] = 0.571
$arr[0,16] = 0
$arr[0,17] = 0
$arr[0,18] = -0.5
$arr[0,19] = 0.45
$ws.Range("K315:AD315").Value = $arr

# Row 316
$ws.Range("B316").Value = 7089400
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 'Odra Opole'
$arr[0,1] = 'Znicz Pruszkw'
$arr[0,2] = 2
$arr[0,3] = 0
$ws.Range("E316:H316").Value = $arr
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = 'H'
$arr[0,1] = 1.615
$arr[0,2] = 3.6
$arr[0,3] = 4.5
$arr[0,4] = 1.65
$arr[0,5] = 3.6
$arr[0,6] = 4.333
$arr[0,7] = -0.75
$arr[0,8] = 1.9
$arr[0,9] = 1.9
$arr[0,10] = 2.25
$arr[0,11] = 1.775
$arr[0,12] = 2.025
$arr[0,13] = 0.6499999999999999
$arr[0,14] = -1
$arr[0,15] = -1
$arr[0,16] = 0.8999999999999999
$arr[0,17] = -1
$arr[0,18] = -0.5
$arr[0,19] = 0.5125
$ws.Range("K316:AD316").Value = $arr

